# Auto-generated edit script: apply numeric updates from scheduled-runner diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 9369.546
$ws.Range("I74").Value = 6702.3076
$ws.Range("J74").Value = 13222.223
$ws.Range("K74").Value = 6702.3076
$ws.Range("L74").Value = 13222.223
$ws.Range("M74").Value = -5766.3076
$ws.Range("N74").Value = -15094.223
$ws.Range("H77").Value = 9369.546
$ws.Range("I77").Value = 6702.3076
$ws.Range("J77").Value = 13222.223
$ws.Range("K77").Value = 33511.538
$ws.Range("L77").Value = 66111.11500000001
$ws.Range("M77").Value = -28831.538
$ws.Range("N77").Value = -75471.11500000001
$ws.Range("H101").Value = 5045.5
$ws.Range("I101").Value = 1330.75
$ws.Range("J101").Value = 9998.5
$ws.Range("K101").Value = 3992.25
$ws.Range("L101").Value = 29995.5
$ws.Range("M101").Value = -2370.25
$ws.Range("N101").Value = -33239.5
$ws.Range("H135").Value = 683.0625
$ws.Range("I135").Value = 710.26666
$ws.Range("K135").Value = 6392.39994
$ws.Range("M135").Value = -3857.39994
$ws.Range("H138").Value = 3002.4092
$ws.Range("I138").Value = 3447.6
$ws.Range("K138").Value = 10342.8
$ws.Range("M138").Value = -5202.799999999999
$ws.Range("H141").Value = 6236.615
$ws.Range("I141").Value = 5634.35
$ws.Range("K141").Value = 16903.05
$ws.Range("M141").Value = -11723.05

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6695.2974
$ws.Range("I32").Value = 6219.971
$ws.Range("K32").Value = 6219.971
$ws.Range("M32").Value = -5932.971
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H110").Value = 5433.0835
$ws.Range("I110").Value = 4645.25
$ws.Range("K110").Value = 4645.25
$ws.Range("M110").Value = -2600.25
$ws.Range("H132").Value = 1379.6383
$ws.Range("I132").Value = 1303.2727
$ws.Range("K132").Value = 3909.8181
$ws.Range("M132").Value = -1379.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 919.9091
$ws.Range("I94").Value = 492
$ws.Range("K94").Value = 492
$ws.Range("M94").Value = -41

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8650.1875
$ws.Range("I58").Value = 2764.8333
$ws.Range("K58").Value = 2764.8333
$ws.Range("M58").Value = -2561.8333
$ws.Range("H62").Value = 4199
$ws.Range("I62").Value = 4665.6665
$ws.Range("K62").Value = 4665.6665
$ws.Range("M62").Value = -4041.6665
$ws.Range("H65").Value = 4199
$ws.Range("I65").Value = 4665.6665
$ws.Range("K65").Value = 23328.3325
$ws.Range("M65").Value = -20208.3325
$ws.Range("H106").Value = 26324.2
$ws.Range("I106").Value = 10000
$ws.Range("J106").Value = 30405.25
$ws.Range("K106").Value = 10000
$ws.Range("L106").Value = 30405.25
$ws.Range("M106").Value = -8738
$ws.Range("N106").Value = -32929.25
$ws.Range("H107").Value = 837.2105
$ws.Range("I107").Value = 877.3077
$ws.Range("J107").Value = 750.3333
$ws.Range("K107").Value = 877.3077
$ws.Range("L107").Value = 750.3333
$ws.Range("M107").Value = 1042.6923
$ws.Range("N107").Value = -4590.3333
$ws.Range("H134").Value = 4717.4443
$ws.Range("I134").Value = 2901.25
$ws.Range("K134").Value = 8703.75
$ws.Range("M134").Value = -6168.75
$ws.Range("H136").Value = 8650.1875
$ws.Range("I136").Value = 2764.8333
$ws.Range("K136").Value = 8294.499899999999
$ws.Range("M136").Value = -5744.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 575
$ws.Range("H98").Value = 530.0769
$ws.Range("I98").Value = 950.6
$ws.Range("K98").Value = 2851.8
$ws.Range("M98").Value = -1353.8
$ws.Range("H131").Value = 12822833
$ws.Range("I131").Value = 83334250
$ws.Range("K131").Value = 250002750
$ws.Range("M131").Value = -249997710
$ws.Range("H141").Value = 4154.7856
$ws.Range("I141").Value = 4320.6924
$ws.Range("J141").Value = 1998
$ws.Range("K141").Value = 12962.0772
$ws.Range("L141").Value = 5994
$ws.Range("M141").Value = -7782.0772
$ws.Range("N141").Value = -16354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 365470.2
$ws.Range("I113").Value = 572296
$ws.Range("K113").Value = 572296
$ws.Range("M113").Value = -570126
$ws.Range("H126").Value = 3734.125
$ws.Range("I126").Value = 2674.7
$ws.Range("K126").Value = 8024.099999999999
$ws.Range("M126").Value = -5554.099999999999
$ws.Range("H132").Value = 6191.5
$ws.Range("I132").Value = 5910.75
$ws.Range("K132").Value = 17732.25
$ws.Range("M132").Value = -15202.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3678
$ws.Range("J16").Value = 4665
$ws.Range("L16").Value = 4665
$ws.Range("N16").Value = -5005
$ws.Range("H55").Value = 50000108
$ws.Range("J55").Value = 134.33333
$ws.Range("L55").Value = 134.33333
$ws.Range("N55").Value = -480.33333
$ws.Range("H93").Value = 25060.727
$ws.Range("I93").Value = 7364.1665
$ws.Range("K93").Value = 7364.1665
$ws.Range("M93").Value = -6116.1665
$ws.Range("H132").Value = 5428
$ws.Range("I132").Value = 5816.32
$ws.Range("J132").Value = 3810
$ws.Range("K132").Value = 17448.96
$ws.Range("L132").Value = 11430
$ws.Range("M132").Value = -14918.96
$ws.Range("N132").Value = -16490

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 50875
$ws.Range("J54").Value = 80872.5
$ws.Range("L54").Value = 80872.5
$ws.Range("N54").Value = -81912.5
$ws.Range("H81").Value = 1096.6666
$ws.Range("I81").Value = 1111.5385
$ws.Range("K81").Value = 2223.077
$ws.Range("M81").Value = -1162.077
$ws.Range("H84").Value = 1096.6666
$ws.Range("I84").Value = 1111.5385
$ws.Range("K84").Value = 11115.385
$ws.Range("M84").Value = -5811.385000000002
$ws.Range("H100").Value = 1419.8572
$ws.Range("I100").Value = 1189.2
$ws.Range("J100").Value = 1996.5
$ws.Range("K100").Value = 2378.4
$ws.Range("L100").Value = 3993
$ws.Range("M100").Value = -1837.4
$ws.Range("N100").Value = -5075
$ws.Range("H132").Value = 5236.8286
$ws.Range("I132").Value = 4593.357
$ws.Range("K132").Value = 13780.071
$ws.Range("M132").Value = -11250.071
$ws.Range("H136").Value = 2431.84
$ws.Range("I136").Value = 2343.0454
$ws.Range("J136").Value = 3083
$ws.Range("K136").Value = 7029.1362
$ws.Range("L136").Value = 9249
$ws.Range("M136").Value = -4479.1362
$ws.Range("N136").Value = -14349
